$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28. This shifts the existing rows 28-49
# down to 29-50 (keeping all their data/formatting intact, mirroring
# Excel's default "insert row, copy formatting from row above" behavior),
# and leaves an empty (but correctly styled) row 28 to be filled in below.
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new weekly record.
$ws.Cells.Item(28, 1).Value = 7
$ws.Cells.Item(28, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(28, 3).Value = "Ñuble"
$ws.Cells.Item(28, 4).Value = 45079
$ws.Cells.Item(28, 5).Value = 16
$ws.Cells.Item(28, 6).Value = 100112043
$ws.Cells.Item(28, 7).Value = "Pepino dulce"
$ws.Cells.Item(28, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 60
$ws.Cells.Item(28, 11).Value = 14000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 14500
$ws.Cells.Item(28, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(28, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(28, 16).Value = 806
$ws.Cells.Item(28, 17).Value = 18
$ws.Cells.Item(28, 18).Value = "Hortaliza"
